# KANSAS_2022 cleanup:
#  1. Rename the header row (A1:D1) from the Spanish descriptive labels to
#     short snake_case field names.
#  2. Normalize Spanish place names in columns A and B to "Proper Case"
#     (Excel's PROPER() function) so connector words like "de", "del", "el",
#     "la", "los", "las", "y" are capitalized like the rest of the words
#     (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"), and fix the
#     stray inner-capital typo "MonteMorelos" -> "Montemorelos".
#  3. Drop the trailing footnote/source rows (663: "Tamano de la muestra..."
#     through 667: "Mayo de 2023") that live below the data table, and
#     shrink the sheet's used range back down to the real data (A1:D661).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Proper-case the state / municipality text columns --------------
$lastDataRow = 661
for ($r = 2; $r -le $lastDataRow; $r++) {
    $stateCell = $ws.Cells.Item($r, 1)
    $stateVal = $stateCell.Value2
    if ($stateVal -ne $null -and $stateVal -ne "") {
        $stateCell.Value = $excel.WorksheetFunction.Proper($stateVal)
    }

    $muniCell = $ws.Cells.Item($r, 2)
    $muniVal = $muniCell.Value2
    if ($muniVal -ne $null -and $muniVal -ne "") {
        $muniCell.Value = $excel.WorksheetFunction.Proper($muniVal)
    }
}

# --- 3. Remove the trailing footnote rows (663-667) below the table ----
$ws.Range("A663:A667").EntireRow.Delete()
